$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 314; this shifts rows 314:361 down to 315:362
# (matching the canonical diff where every existing record from row 314
# onward moves down by one row, and a brand-new weekly price record is
# inserted at row 314).
$ws.Rows("314:314").Insert()

# Populate the newly inserted row 314 with the new record's data. The
# "template" columns (A,B,C,E,F,G,H,I,J,K,R) repeat the same market /
# product / variety / origin values used by all the surrounding Hayward
# Kiwi - Vega Modelo de Temuco rows.
$ws.Range("A314").Value = 10
$ws.Range("B314").Value = "Vega Modelo de Temuco"
$ws.Range("C314").Value = "La Araucanía"
$ws.Range("D314").Value = 44505
$ws.Range("E314").Value = 9
$ws.Range("F314").Value = "Fruta"
$ws.Range("G314").Value = 100101
$ws.Range("H314").Value = "Berries"
$ws.Range("I314").Value = 100101007
$ws.Range("J314").Value = "Kiwi"
$ws.Range("K314").Value = "Hayward"
$ws.Range("L314").Value = "Especial"
$ws.Range("M314").Value = 55
$ws.Range("N314").Value = 20000
$ws.Range("O314").Value = 20000
$ws.Range("P314").Value = 20000
$ws.Range("Q314").Value = "$/caja 15 kilos"
$ws.Range("R314").Value = "Región de O'Higgins"
$ws.Range("S314").Value = 1333
$ws.Range("T314").Value = 15

# Match the date-formatted style already used by the rest of column D.
$ws.Range("D314").NumberFormat = $ws.Range("D315").NumberFormat
